$wb = $excel.ActiveWorkbook
$wsGlavnye = $wb.Worksheets.Item(2)   # "Главные" (main referees)
$wsLineynye = $wb.Worksheets.Item(3)  # "Линейные" (linesmen)

# --- Updated statistic values (rows with new games/PIM data) ---
# $wsGlavnye row 4
$wsGlavnye.Range("C4").Value = 17
$wsGlavnye.Range("D4").Value = 288
$wsGlavnye.Range("E4").Value = 125
$wsGlavnye.Range("F4").Value = 163
$wsGlavnye.Range("G4").Value = 16.94
$wsGlavnye.Range("H4").Value = 7.35
$wsGlavnye.Range("I4").Value = 9.59
$wsGlavnye.Range("J4").Value = 60
$wsGlavnye.Range("K4").Value = 69
$wsGlavnye.Range("W4").Value = 12
# $wsGlavnye row 5
$wsGlavnye.Range("C5").Value = 23
$wsGlavnye.Range("D5").Value = 367
$wsGlavnye.Range("E5").Value = 199
$wsGlavnye.Range("F5").Value = 168
$wsGlavnye.Range("G5").Value = 15.96
$wsGlavnye.Range("H5").Value = 8.65
$wsGlavnye.Range("I5").Value = 7.3
$wsGlavnye.Range("J5").Value = 97
$wsGlavnye.Range("K5").Value = 84
# $wsGlavnye row 9
$wsGlavnye.Range("C9").Value = 23
$wsGlavnye.Range("D9").Value = 366
$wsGlavnye.Range("E9").Value = 197
$wsGlavnye.Range("F9").Value = 169
$wsGlavnye.Range("G9").Value = 15.91
$wsGlavnye.Range("H9").Value = 8.57
$wsGlavnye.Range("I9").Value = 7.35
$wsGlavnye.Range("J9").Value = 96
$wsGlavnye.Range("K9").Value = 82
$wsGlavnye.Range("W9").Value = 20
# $wsGlavnye row 10
$wsGlavnye.Range("C10").Value = 15
$wsGlavnye.Range("D10").Value = 276
$wsGlavnye.Range("E10").Value = 138
$wsGlavnye.Range("F10").Value = 138
$wsGlavnye.Range("G10").Value = 18.4
$wsGlavnye.Range("H10").Value = 9.199999999999999
$wsGlavnye.Range("I10").Value = 9.199999999999999
$wsGlavnye.Range("J10").Value = 69
$wsGlavnye.Range("K10").Value = 59
# $wsGlavnye row 11
$wsGlavnye.Range("C11").Value = 16
$wsGlavnye.Range("D11").Value = 423
$wsGlavnye.Range("E11").Value = 198
$wsGlavnye.Range("F11").Value = 225
$wsGlavnye.Range("G11").Value = 26.44
$wsGlavnye.Range("H11").Value = 12.38
$wsGlavnye.Range("I11").Value = 14.06
$wsGlavnye.Range("J11").Value = 89
$wsGlavnye.Range("K11").Value = 75
# $wsGlavnye row 14
$wsGlavnye.Range("C14").Value = 14
$wsGlavnye.Range("D14").Value = 171
$wsGlavnye.Range("E14").Value = 90
$wsGlavnye.Range("G14").Value = 12.21
$wsGlavnye.Range("H14").Value = 6.43
$wsGlavnye.Range("I14").Value = 5.79
$wsGlavnye.Range("J14").Value = 45
# $wsGlavnye row 15
$wsGlavnye.Range("C15").Value = 15
$wsGlavnye.Range("D15").Value = 206
$wsGlavnye.Range("F15").Value = 116
$wsGlavnye.Range("G15").Value = 13.73
$wsGlavnye.Range("H15").Value = 6
$wsGlavnye.Range("I15").Value = 7.73
$wsGlavnye.Range("K15").Value = 58
# $wsGlavnye row 23
$wsGlavnye.Range("C23").Value = 14
$wsGlavnye.Range("D23").Value = 181
$wsGlavnye.Range("F23").Value = 119
$wsGlavnye.Range("G23").Value = 12.93
$wsGlavnye.Range("H23").Value = 4.43
$wsGlavnye.Range("I23").Value = 8.5
$wsGlavnye.Range("K23").Value = 47
# $wsLineynye row 3
$wsLineynye.Range("C3").Value = 21
$wsLineynye.Range("D3").Value = 308
$wsLineynye.Range("E3").Value = 146
$wsLineynye.Range("F3").Value = 162
$wsLineynye.Range("G3").Value = 14.67
$wsLineynye.Range("H3").Value = 6.95
$wsLineynye.Range("I3").Value = 7.71
$wsLineynye.Range("J3").Value = 73
$wsLineynye.Range("K3").Value = 66
$wsLineynye.Range("W3").Value = 10
# $wsLineynye row 6
$wsLineynye.Range("C6").Value = 13
$wsLineynye.Range("D6").Value = 241
$wsLineynye.Range("E6").Value = 100
$wsLineynye.Range("F6").Value = 141
$wsLineynye.Range("G6").Value = 18.54
$wsLineynye.Range("H6").Value = 7.69
$wsLineynye.Range("I6").Value = 10.85
$wsLineynye.Range("J6").Value = 50
$wsLineynye.Range("K6").Value = 68
# $wsLineynye row 7
$wsLineynye.Range("C7").Value = 13
$wsLineynye.Range("D7").Value = 215
$wsLineynye.Range("F7").Value = 143
$wsLineynye.Range("G7").Value = 16.54
$wsLineynye.Range("H7").Value = 5.54
$wsLineynye.Range("I7").Value = 11
$wsLineynye.Range("K7").Value = 44
# $wsLineynye row 12
$wsLineynye.Range("C12").Value = 20
$wsLineynye.Range("D12").Value = 348
$wsLineynye.Range("E12").Value = 167
$wsLineynye.Range("F12").Value = 181
$wsLineynye.Range("G12").Value = 17.4
$wsLineynye.Range("H12").Value = 8.35
$wsLineynye.Range("I12").Value = 9.050000000000001
$wsLineynye.Range("J12").Value = 76
$wsLineynye.Range("K12").Value = 83
# $wsLineynye row 14
$wsLineynye.Range("C14").Value = 23
$wsLineynye.Range("D14").Value = 382
$wsLineynye.Range("E14").Value = 192
$wsLineynye.Range("G14").Value = 16.61
$wsLineynye.Range("H14").Value = 8.35
$wsLineynye.Range("I14").Value = 8.26
$wsLineynye.Range("J14").Value = 96
# $wsLineynye row 19
$wsLineynye.Range("C19").Value = 21
$wsLineynye.Range("D19").Value = 377
$wsLineynye.Range("E19").Value = 176
$wsLineynye.Range("F19").Value = 201
$wsLineynye.Range("H19").Value = 8.380000000000001
$wsLineynye.Range("I19").Value = 9.57
$wsLineynye.Range("J19").Value = 83
$wsLineynye.Range("K19").Value = 88
# $wsLineynye row 20
$wsLineynye.Range("C20").Value = 16
$wsLineynye.Range("D20").Value = 264
$wsLineynye.Range("E20").Value = 135
$wsLineynye.Range("G20").Value = 16.5
$wsLineynye.Range("H20").Value = 8.44
$wsLineynye.Range("I20").Value = 8.06
$wsLineynye.Range("J20").Value = 65
# $wsLineynye row 21
$wsLineynye.Range("C21").Value = 25
$wsLineynye.Range("D21").Value = 517
$wsLineynye.Range("F21").Value = 304
$wsLineynye.Range("G21").Value = 20.68
$wsLineynye.Range("H21").Value = 8.52
$wsLineynye.Range("I21").Value = 12.16
$wsLineynye.Range("K21").Value = 122

# --- Refresh "as_of_utc" timestamp (column AA) for every data row (2-26) on both sheets ---
$newTimestamp = "2025-11-09 03:02:37"
for ($r = 2; $r -le 26; $r++) {
    $wsGlavnye.Range("AA" + $r).Value = $newTimestamp
    $wsLineynye.Range("AA" + $r).Value = $newTimestamp
}
